# Add six new cafeteria rows (20-25) to Sheet1, matching the upstream
# "Add files via upload" commit that appended more cafeteria_name rows
# to the Cafeteria.xlsx data sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pre-format the new range the same way the existing data rows are
# formatted (centered alignment == cellXfs index 1 in the original file)
# so the new cells reuse the existing style instead of minting a new one.
$ws.Range("A22:C27").HorizontalAlignment = -4108  # xlCenter

# Populate the shared-string column (B) first, in the same order the
# strings were originally authored, so the shared-string table indices
# line up with the upstream workbook (22:西园食堂 .. 27:榕园食堂).
$ws.Range("B22").Value = "西园食堂"
$ws.Range("B23").Value = "东园食堂"
$ws.Range("B25").Value = "荔园食堂"
$ws.Range("B26").Value = "槿园食堂"
$ws.Range("B27").Value = "若海食堂"
$ws.Range("B24").Value = "榕园食堂"

# cafeteria_id (A) and campus_id (C) columns.
$ws.Range("A22").Value = 20
$ws.Range("C22").Value = 3

$ws.Range("A23").Value = 21
$ws.Range("C23").Value = 3

$ws.Range("A24").Value = 22
$ws.Range("C24").Value = 4

$ws.Range("A25").Value = 23
$ws.Range("C25").Value = 4

$ws.Range("A26").Value = 24
$ws.Range("C26").Value = 4

$ws.Range("A27").Value = 25
$ws.Range("C27").Value = 4

# Match the saved selection/scroll state recorded in the upstream diff.
$ws.Range("B24").Select()
